$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new version row (row 4) with the 0.1.1 release info.
$ws.Range("A4").Value = "0.1.1"
$ws.Range("B4").Value = $ws.Range("B3").Value2
$ws.Range("D4").Value = " -Changed the sorting order for the Group selection to be proportional from the Direct Fitness/Indirect Fitness`n-Implemented  group selection correctly."
$ws.Range("C4").Value = "`n-Agrupation and desagrupation to be done in functions.`n-Change reproduction and distribution to two parts.`n-UI: Delete rows according to working functionality.`n-UI: condicionate IF to be associated.`n-Implement mutations.`n-Implement save and load configurations.`n-Document every function.`n-The program does not work with 2 or less niches.`n-When an actor dies, the recipient should deassociate.`n"
$ws.Range("E4").Value = $ws.Range("E3").Value2
$ws.Range("F4").Value = $ws.Range("F3").Value2
$ws.Range("G4").Value = $ws.Range("G3").Value2

# Copy the formatting that the banding row (row 2) uses for the new row.
$ws.Range("A2:G2").Copy()
$ws.Range("A4:G4").PasteSpecial(-4122)

# Turn on wrap text for the whole table (columns A-G) like Excel does when
# the user selects the columns and toggles Wrap Text.
$ws.Columns("A:G").WrapText = $true

# Let Excel auto-fit the new row's height instead of keeping a custom one.
$ws.Rows("4:4").AutoFit()

# Update selection / frozen pane state to match final view.
$ws.Range("C4").Select()
